$d = $word.ActiveDocument

# The "FORM KELENGKAPAN ADMINISTRASI" checklist table is the 4th table in
# the document. Row 13 is the "11 - Penambahan Peta Titk Pengelolaan dan
# Titik Pemantauan" verification row, immediately after the "Persetujuan
# Teknis" (row 10) entry and before the "Catatan untuk perbaikan
# pemrakarsa" row. Locate it defensively by its first-cell marker text
# rather than a hard-coded index, then delete the whole row.

$tbl = $d.Tables.Item(4)

for ($r = $tbl.Rows.Count; $r -ge 1; $r--) {
    $row = $tbl.Rows.Item($r)
    $marker = $row.Cells.Item(2).Range.Text
    if ($marker -like "*Penambahan Peta Titk Pengelolaan dan Titik Pemantauan*") {
        $row.Delete()
    }
}
